$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 102
$ws.Range("B102").Value = 5461531
$ws.Range("F102").Value = "Vasas SC"
$ws.Range("G102").Value = "MOL Fehervar FC"
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = "D"
$ws.Range("K102").Value = 6
$ws.Range("L102").Value = 4.5
$ws.Range("M102").Value = 1.4
$ws.Range("N102").Value = 5.75
$ws.Range("O102").Value = 4.75
$ws.Range("P102").Value = 1.5
$ws.Range("Q102").Value = 1
$ws.Range("R102").Value = 2.05
$ws.Range("S102").Value = 1.8
$ws.Range("U102").Value = 1.825
$ws.Range("V102").Value = 2.025
$ws.Range("W102").Value = -1
$ws.Range("X102").Value = 3.75
$ws.Range("Z102").Value = 1.05
$ws.Range("AC102").Value = 1.025

# Row 103
$ws.Range("B103").Value = 5470380
$ws.Range("F103").Value = "Puskas Academy"
$ws.Range("G103").Value = "Budapest Honved"
$ws.Range("H103").Value = 2
$ws.Range("I103").Value = 1
$ws.Range("J103").Value = "H"
$ws.Range("K103").Value = 2
$ws.Range("L103").Value = 3.6
$ws.Range("M103").Value = 3.6
$ws.Range("N103").Value = 1.75
$ws.Range("O103").Value = 3.8
$ws.Range("P103").Value = 4.5
$ws.Range("Q103").Value = -0.75
$ws.Range("R103").Value = 2
$ws.Range("S103").Value = 1.85
$ws.Range("U103").Value = 2
$ws.Range("V103").Value = 1.85
$ws.Range("W103").Value = 0.75
$ws.Range("X103").Value = -1
$ws.Range("Z103").Value = 0.5
$ws.Range("AA103").Value = -0.5
$ws.Range("AB103").Value = 0.5
$ws.Range("AC103").Value = -0.5

# Row 104
$ws.Range("B104").Value = 5461530
$ws.Range("F104").Value = "Debreceni VSC"
$ws.Range("G104").Value = "Ujpest"
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 1.727
$ws.Range("L104").Value = 3.5
$ws.Range("M104").Value = 4
$ws.Range("N104").Value = 1.7
$ws.Range("O104").Value = 4
$ws.Range("P104").Value = 4.75
$ws.Range("R104").Value = 1.875
$ws.Range("S104").Value = 1.975
$ws.Range("U104").Value = 1.975
$ws.Range("V104").Value = 1.875
$ws.Range("W104").Value = 0.7
$ws.Range("Z104").Value = 0.875
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = -1
$ws.Range("AC104").Value = 0.875

# Row 232
$ws.Range("N232").Value = 2.05
$ws.Range("O232").Value = 3.6
$ws.Range("P232").Value = 3.1
$ws.Range("Q232").Value = -0.25
$ws.Range("R232").Value = 1.85
$ws.Range("S232").Value = 2
$ws.Range("U232").Value = 1.925
$ws.Range("V232").Value = 1.925

# Row 235
$ws.Range("R235").Value = 1.85
$ws.Range("S235").Value = 2
$ws.Range("U235").Value = 1.825
$ws.Range("V235").Value = 2.025
